# Autorizzazione_Paesaggistica_146.docx — "Modifiche modelli di stampa"
#
# The title line on the second page currently reads:
#     n.          del
# (a bare placeholder with nothing but blank space for the act number
# and the release date). This adds the two mail-merge style tokens
# "[numero_titolo]" and "[data_rilascio_titolo]" into that blank space,
# keeping the original amount of surrounding whitespace intact, and
# drops a hidden "_GoBack" bookmark right after the new date token
# (immediately before the closing "]"), matching where Word leaves the
# edit-position bookmark after the last text insertion in the document.

$d = $word.ActiveDocument

$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "n.*del *") {
        $titlePara = $p
        break
    }
}

if ($titlePara -eq $null) {
    throw "Could not locate the 'n. ... del ...' title paragraph"
}

$range = $titlePara.Range
$replaced = $range.Find.Execute(
    "n.          del ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "n.    [numero_titolo]      del [data_rilascio_titolo]",
    2)

if (-not $replaced) {
    throw "Find/Replace for the title placeholder text did not match"
}

# Re-find the freshly inserted date token so we can drop the _GoBack
# bookmark right before its closing bracket, exactly as Word would
# leave the caret/bookmark after typing the text in interactively.
$tokenRange = $d.Content
$tokenRange.Find.Execute("[data_rilascio_titolo]") | Out-Null
$bookmarkSpot = $d.Range($tokenRange.End - 1, $tokenRange.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot) | Out-Null
